$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Swap the OS / EFS header labels in row 1
$ws.Range("B1").Value = "Significant CpG probes in EFS"
$ws.Range("C1").Value = "Significant CpG probes in OS"
# D1 "Significant overlapping CpG probes" stays the same

# Update row 2 label and values
$ws.Range("A2").Value = "Updated Risk Group"
$ws.Range("B2").Value = 112
$ws.Range("C2").Value = 189
$ws.Range("D2").Value = 17

# Remove row 3 entirely (was "Adjusted" / 167 / 55 / 17)
$ws.Rows("3:3").Delete()
